$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: replace the old INDEX(A1:B5,1,2) text formula with a formula that
# errors out with #N/A (INDEX given an invalid reference + a volatile SHEET call).
$ws.Range("D6").Formula = "=INDEX(9, _xlfn.SHEET(2))"

# Rows 7 and 8 (the "2+3" and "3+2" helper formulas) are removed entirely.
$ws.Range("D7:D8").ClearContents()
